# Reorder slides: move the "Context diagram" slide (originally slide 8)
# to position 2 (right after the title slide), and delete the
# "User manual" slide (originally slide 9, the last slide).
$p = $ppt.ActivePresentation

# Move slide 8 ("Context diagram") to become the 2nd slide.
$contextDiagram = $p.Slides.Item(8)
$contextDiagram.MoveTo(2)

# After the move, the old "User manual" slide is still the last slide
# (it was slide 9 and nothing after it moved), so it's still at index 9.
$userManual = $p.Slides.Item(9)
$userManual.Delete()
